# Applies the "Change colors, add remaining cards" edit to the allies sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns for the split-out attribute icon/sign/number fields.
$ws.Range("I1").Value = "Icon"
$ws.Range("J1").Value = "Sign"
$ws.Range("K1").Value = "Number"
$ws.Range("L1").Value = "Icon_2"
$ws.Range("M1").Value = "Sign_2"
$ws.Range("N1").Value = "Number_2"

# Row 14 - Rabble-rouser: Income - 2  =>  Single_Attribute / money / - / 2
$ws.Range("H14").Value = "Single_Attribute"
$ws.Range("I14").Value = "money"
$ws.Range("J14").Value = "-"
$ws.Range("K14").Value = 2

# Row 15 - Researcher: Peek - 1  =>  Single_Attribute / eye / + / 1
$ws.Range("H15").Value = "Single_Attribute"
$ws.Range("I15").Value = "eye"
$ws.Range("J15").Value = "+"
$ws.Range("K15").Value = 1

# Row 16 - Tax Collector: Income + 1  =>  Single_Attribute / money / + / 1
$ws.Range("H16").Value = "Single_Attribute"
$ws.Range("I16").Value = "money"
$ws.Range("J16").Value = "+"
$ws.Range("K16").Value = 1

# Row 17 - Ambassador: Trade - 1  =>  Single_Attribute / trade / + / 1
$ws.Range("H17").Value = "Single_Attribute"
$ws.Range("I17").Value = "trade"
$ws.Range("J17").Value = "+"
$ws.Range("K17").Value = 1

# Row 18 - Philanthropist: Income + 1, Peek - 1  =>  Double_Attribute / money + 1 / eye - 1
$ws.Range("H18").Value = "Double_Attribute"
$ws.Range("I18").Value = "money"
$ws.Range("J18").Value = "+"
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = "eye"
$ws.Range("M18").Value = "-"
$ws.Range("N18").Value = 1

# Row 19 - Merchant: Income + 1, Trade - 1  =>  Double_Attribute / money + 1 / trade - 1
$ws.Range("H19").Value = "Double_Attribute"
$ws.Range("I19").Value = "money"
$ws.Range("J19").Value = "+"
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = "trade"
$ws.Range("M19").Value = "-"
$ws.Range("N19").Value = 1

# Match the saved selection state from the edited workbook.
$ws.Range("I20").Select()
